$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (357) down to the new rows (358-366)
$ws.Range("A357").Copy()
$ws.Range("A358:A366").PasteSpecial(-4122)

# Fill in the new data rows
$ws.Cells.Item(358, 1).Value = 44432
$ws.Cells.Item(358, 2).Value = 1
$ws.Cells.Item(358, 3).Value = 21
$ws.Cells.Item(358, 4).Value = 61.10690799045568

$ws.Cells.Item(359, 1).Value = 44433
$ws.Cells.Item(359, 2).Value = 0
$ws.Cells.Item(359, 3).Value = 20
$ws.Cells.Item(359, 4).Value = 58.19705522900541

$ws.Cells.Item(360, 1).Value = 44434
$ws.Cells.Item(360, 2).Value = 8
$ws.Cells.Item(360, 3).Value = 23
$ws.Cells.Item(360, 4).Value = 66.92661351335623

$ws.Cells.Item(361, 1).Value = 44435
$ws.Cells.Item(361, 2).Value = 4
$ws.Cells.Item(361, 3).Value = 20
$ws.Cells.Item(361, 4).Value = 58.19705522900541

$ws.Cells.Item(362, 1).Value = 44436
$ws.Cells.Item(362, 2).Value = 0
$ws.Cells.Item(362, 3).Value = 20
$ws.Cells.Item(362, 4).Value = 58.19705522900541

$ws.Cells.Item(363, 1).Value = 44437
$ws.Cells.Item(363, 2).Value = 3
$ws.Cells.Item(363, 3).Value = 18
$ws.Cells.Item(363, 4).Value = 52.37734970610487

$ws.Cells.Item(364, 1).Value = 44438
$ws.Cells.Item(364, 2).Value = 0
$ws.Cells.Item(364, 3).Value = 16
$ws.Cells.Item(364, 4).Value = 46.55764418320432

$ws.Cells.Item(365, 1).Value = 44439
$ws.Cells.Item(365, 2).Value = 0
$ws.Cells.Item(365, 3).Value = 15
$ws.Cells.Item(365, 4).Value = 43.64779142175406

$ws.Cells.Item(366, 1).Value = 44440
$ws.Cells.Item(366, 2).Value = 0
$ws.Cells.Item(366, 3).Value = 15
$ws.Cells.Item(366, 4).Value = 43.64779142175406
